$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New descriptive stats values for columns B (F1-Score), C (AUC), D (Acurácia)
# keyed by row number
$updates = @{
    2  = @{ B = 0.6667; C = 0.7375; D = 0.7692 }
    3  = @{ B = 0.75;   C = 0.8;    D = 0.8462 }
    4  = @{ B = 0.6667; C = 0.7375; D = 0.7692 }
    6  = @{ B = 0.2857; C = 0.5833; D = 0.6154 }
    7  = @{ B = 0.4;    C = 0.5238; D = 0.5385 }
    8  = @{ B = 0.5714; C = 0.7;    D = 0.75 }
    10 = @{ B = 0.5714; C = 0.7;    D = 0.75 }
    11 = @{ B = 0.8889; C = 0.9286; D = 0.9091 }
    12 = @{ B = 0.8889; C = 0.9286; D = 0.9091 }
    13 = @{ B = 0.8;    C = 0.8571; D = 0.8182 }
    14 = @{ B = 0.6;    C = 0.6333; D = 0.6364 }
    15 = @{ B = 0.5;    C = 0.6167; D = 0.6364 }
    16 = @{ B = 0.4444; C = 0.5333; D = 0.5455 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
}
